$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column N (22-jun), matching the text style of the other
# date headers in row 1 (column M = 21-jun).
$ws.Range("N1").Value = "22-jun"
$ws.Range("N1").NumberFormat = $ws.Range("M1").NumberFormat

# New numeric values for column N, rows 2-18.
$values = @{
    2  = 0
    3  = 12.696568499229443
    4  = 16.021400258648665
    5  = 15.505918927536154
    6  = 0
    7  = 8.6621980147994719
    8  = 7.6880963776035349
    9  = 15.641182915824235
    10 = 11.262975680222553
    11 = 11.072010352859856
    12 = 0
    13 = 13.502219378486847
    14 = 0
    15 = 0
    16 = 15.33349861184303
    17 = 0
    18 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 14).Value = $values[$row]
}

# Match the selection left after the edit.
$ws.Range("P8").Select() | Out-Null
